$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/date range) ---
$ws.Range("A8").Value = "Volume 32   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/29/2025  Through  10/5/2025"

# --- Donor cells used as format/value templates for text<->number conversions ---
# C14 = text "0" (style 13, shared string 20)
# E14 = text "***.*" (style 13, shared string 21)
# I14 = number style 14 (#,##0)
# L14 = number style 15 (decimal)

$ws.Range("N14").Value = -95.454545454545
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 28
$ws.Range("K15").Value = 3.703703703703
$ws.Range("L15").Value = 12
$ws.Range("M15").Value = 133.333333333333
$ws.Range("N15").Value = -22.222222222222
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -87.5
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = -34.782608695652
$ws.Range("I16").Value = 146
$ws.Range("J16").Value = 236
$ws.Range("K16").Value = -38.135593220339
$ws.Range("L16").Value = -30.143540669856
$ws.Range("M16").Value = -41.129032258064
$ws.Range("N16").Value = -89.022556390977
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -41.666666666666
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = -9.756097560975
$ws.Range("I17").Value = 434
$ws.Range("J17").Value = 448
$ws.Range("K17").Value = -3.125
$ws.Range("L17").Value = 12.435233160621
$ws.Range("M17").Value = 101.860465116279
$ws.Range("N17").Value = -27.906976744186
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 15
$ws.Range("H18").Value = 15.384615384615
$ws.Range("I18").Value = 170
$ws.Range("J18").Value = 153
$ws.Range("K18").Value = 11.111111111111
$ws.Range("L18").Value = -0.584795321637
$ws.Range("M18").Value = -54.301075268817
$ws.Range("N18").Value = -90.122022080185
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 15.789473684210
$ws.Range("F19").Value = 85
$ws.Range("G19").Value = 91
$ws.Range("H19").Value = -6.593406593406
$ws.Range("I19").Value = 709
$ws.Range("J19").Value = 773
$ws.Range("K19").Value = -8.279430789133
$ws.Range("L19").Value = 0.997150997150
$ws.Range("M19").Value = 41.8
$ws.Range("N19").Value = -16.489988221437
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 36
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = 56.521739130434
$ws.Range("I20").Value = 271
$ws.Range("J20").Value = 274
$ws.Range("K20").Value = -1.094890510948
$ws.Range("L20").Value = -19.822485207100
$ws.Range("M20").Value = 25.462962962963
$ws.Range("N20").Value = -86.863790596219
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 46
$ws.Range("E21").Value = -13.043478260869
$ws.Range("F21").Value = 189
$ws.Range("G21").Value = 192
$ws.Range("H21").Value = -1.5625
$ws.Range("I21").Value = 1759
$ws.Range("J21").Value = 1911
$ws.Range("K21").Value = -7.953950811093
$ws.Range("L21").Value = -4.193899782135
$ws.Range("M21").Value = 12.038216560509
$ws.Range("N21").Value = -73.441038804167
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("J22").Value = 35
$ws.Range("K22").Value = -17.142857142857
$ws.Range("L22").Value = -12.121212121212
$ws.Range("M22").Value = 20.833333333333
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 150
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 8.333333333333
$ws.Range("I23").Value = 161
$ws.Range("J23").Value = 168
$ws.Range("K23").Value = -4.166666666666
$ws.Range("L23").Value = -18.274111675126
$ws.Range("M23").Value = 37.606837606837
$ws.Range("C24").Value = 61
$ws.Range("D24").Value = 66
$ws.Range("E24").Value = -7.575757575757
$ws.Range("F24").Value = 188
$ws.Range("G24").Value = 222
$ws.Range("H24").Value = -15.315315315315
$ws.Range("I24").Value = 1717
$ws.Range("J24").Value = 1765
$ws.Range("K24").Value = -2.719546742209
$ws.Range("L24").Value = -5.242825607064
$ws.Range("M24").Value = 47.889750215331
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = -9.375
$ws.Range("F25").Value = 91
$ws.Range("G25").Value = 121
$ws.Range("H25").Value = -24.793388429752
$ws.Range("I25").Value = 857
$ws.Range("J25").Value = 904
$ws.Range("K25").Value = -5.199115044247
$ws.Range("L25").Value = 3.878787878787
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 43.75
$ws.Range("F26").Value = 81
$ws.Range("G26").Value = 62
$ws.Range("H26").Value = 30.645161290322
$ws.Range("I26").Value = 673
$ws.Range("J26").Value = 722
$ws.Range("K26").Value = -6.786703601108
$ws.Range("L26").Value = 3.062787136294
$ws.Range("M26").Value = -1.464128843338
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4163)
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("I27").Value = 34
$ws.Range("J27").Value = 41
$ws.Range("K27").Value = -17.073170731707
$ws.Range("L27").Value = -12.820512820512
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -25
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -54.545454545454
$ws.Range("I28").Value = 71
$ws.Range("J28").Value = 74
$ws.Range("K28").Value = -4.054054054054
$ws.Range("L28").Value = -21.111111111111
$ws.Range("C14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G29").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H29").PasteSpecial(-4163)
$ws.Range("L29").Value = -60.869565217391
$ws.Range("N29").Value = -73.529411764705
$ws.Range("C14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G30").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("L30").Value = -57.142857142857
$ws.Range("N30").Value = -71.875
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F31").PasteSpecial(-4163)
$ws.Range("H31").Value = -100
$ws.Range("C14").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F33").PasteSpecial(-4163)
$ws.Range("I14").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("G33").Value = 1
$ws.Range("L14").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").Value = -100
$ws.Range("J33").Value = 4
$ws.Range("K33").Value = 75

$excel.CutCopyMode = 0

